$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - update existing headers and add new ones through column M
$headers = @("HTNO", "Edge Analytics", "Cyber Security", "Machine Learning", "Software Project Management", "Human Computer Interaction", "Renewable Energy Source", "Edge Analytics Lab", "Machine Learning Lab", "Internship / Mini Project", "Total", "SGPA", "CGPA")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows - all values are text (hall-ticket numbers, marks, totals, SGPA), so force
# the cells to text format before assigning to avoid Excel auto-converting numeric-looking
# strings into numeric values.
$data = @(
    @("20J41A6901", "87", "80", "84", "74", "80", "87", "99", "100", "99", "790", "9.08"),
    @("20J41A6902", "80", "77", "72", "71", "79", "77", "94", "98", "98", "746", "8.58")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 2, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c]
        $cell.ClearFormats()
    }
}
